# edit.ps1
# Applies the commit "add translations core, damageFlashComponent, particles of shot"
# to dev_notes.xlsx:
#   - TODO Before 0.0.1: mark row 11 (C11) status "todo" -> "in-progress"
#   - Logs: append 4 new log rows (28-31) describing the work
#   - Selection / active-tab bookkeeping: Logs becomes the active sheet/tab

$wb = $excel.ActiveWorkbook

$wsTodo = $wb.Worksheets.Item("TODO Before 0.0.1")
$wsLogs = $wb.Worksheets.Item("Logs")

# --- TODO Before 0.0.1: update status of row 11 from "todo" to "in-progress"
$wsTodo.Range("C11").Value = "in-progress"

# --- Logs: add the four new dated entries right after the existing last row (27)
$wsLogs.Range("A27").Copy($wsLogs.Range("A28:A31"))

$wsLogs.Range("A28").Value = 45445
$wsLogs.Range("B28").Value = "work on particles for projectiles - looks good but need finilize for simple attack "

$wsLogs.Range("A29").Value = 45451
$wsLogs.Range("B29").Value = "work on partickles, explosion particles, animation of damage for player and enemies, add lighthing, shaders and fog "

$wsLogs.Range("A30").Value = 45452
$wsLogs.Range("B30").Value = "start work on damage flash component"

$wsLogs.Range("A31").Value = 45455
$wsLogs.Range("B31").Value = "Implement damage flash component - for enemies and player, improve it"

# --- View / selection bookkeeping
# TODO Before 0.0.1 keeps its own scroll/selection state, but is no longer the
# tab shown when the workbook is (re)opened.
$wsTodo.Activate()
$wsTodo.Range("C11").Select()

# Logs becomes the active sheet/tab, with its selection on the new last entry.
$wsLogs.Activate()
$wsLogs.Range("B32").Select()
